# Add a new time-entry row (row 13) to Sheet1, extend the two SUM formulas
# that roll up column A to include it, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data point appended below the existing entries (row 12 -> row 13).
# Match the existing number format used by A3:A12 (h:mm:ss time format).
$ws.Range("A13").Value = 0.032361111111111111
$ws.Range("A13").NumberFormat = $ws.Range("A12").NumberFormat

# Extend the rollup formulas so they include the newly added row.
$ws.Range("C2").Formula = "=SUM(A2:A13)"
$ws.Range("B3").Formula = "=SUM(A9:A13)"

# Move the active selection to match the edited workbook's saved state.
$ws.Range("H8").Select()

$wb.Save()
